$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 8-12 (data shrinks from 12 games to 7 games)
$ws.Range("A8:C12").EntireRow.Delete()

# Update header row
$ws.Range("A1").Value() = "NBA, Wednesday 28th Feb 2024"
$ws.Range("B1").Value() = "Ballgorithm"
$ws.Range("C1").Value() = "ESPN"

# Update the matchup/prediction rows with the new day's games
$ws.Range("A2").Value() = "New Orleans Pelicans (35-24) vs Indiana Pacers (33-27)"
$ws.Range("B2").Value() = "Indiana Pacers (61.29%)"
$ws.Range("C2").Value() = "Indiana Pacers (58.8%)"

$ws.Range("A3").Value() = "Memphis Grizzlies (20-38) vs Minnesota Timberwolves (41-17)"
$ws.Range("B3").Value() = "Minnesota Timberwolves (77.78%)"
$ws.Range("C3").Value() = "Minnesota Timberwolves (88.9%)"

$ws.Range("A4").Value() = "Dallas Mavericks (33-25) vs Toronto Raptors (22-36)"
$ws.Range("B4").Value() = "Dallas Mavericks (58.06%)"
$ws.Range("C4").Value() = "Dallas Mavericks (57.2%)"

$ws.Range("A5").Value() = "Cleveland Cavaliers (38-19) vs Chicago Bulls (27-31)"
$ws.Range("B5").Value() = "Cleveland Cavaliers (66.67%)"
$ws.Range("C5").Value() = "Cleveland Cavaliers (73.0%)"

$ws.Range("A6").Value() = "Sacramento Kings (33-24) vs Denver Nuggets (39-19)"
$ws.Range("B6").Value() = "Denver Nuggets (81.48%)"
$ws.Range("C6").Value() = "Denver Nuggets (70.5%)"

$ws.Range("A7").Value() = "Los Angeles Lakers (32-28) vs Los Angeles Clippers (37-19)"
$ws.Range("B7").Value() = "Los Angeles Clippers (74.07%)"
$ws.Range("C7").Value() = "Los Angeles Clippers (74.8%)"
